$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057552590131079
$ws.Range("D2").Value = 1.060863525316217
$ws.Range("E2").Value = 1.053643488690419
$ws.Range("F2").Value = 1.069306891305746
$ws.Range("I2").Value = 1.051976358614207
$ws.Range("J2").Value = 1.06254841701045
$ws.Range("K2").Value = 1.063589147592431
$ws.Range("L2").Value = 1.056388890393098
$ws.Range("M2").Value = 1.072009749352477
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058801790855001
$ws.Range("D3").Value = 1.061866529510889
$ws.Range("E3").Value = 1.054722483998048
$ws.Range("F3").Value = 1.070507465402748
$ws.Range("I3").Value = 1.052416036808772
$ws.Range("J3").Value = 1.063448691681089
$ws.Range("K3").Value = 1.064406184900866
$ws.Range("L3").Value = 1.057280290439264
$ws.Range("M3").Value = 1.07302551511633
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059609601515911
$ws.Range("D4").Value = 1.062515043672138
$ws.Range("E4").Value = 1.05542034430982
$ws.Range("F4").Value = 1.071284212715246
$ws.Range("I4").Value = 1.052698965749616
$ws.Range("J4").Value = 1.06403019862639
$ws.Range("K4").Value = 1.064933753438698
$ws.Range("L4").Value = 1.057856169056119
$ws.Range("M4").Value = 1.073682104331751
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059949086695822
$ws.Range("D5").Value = 1.062787561587332
$ws.Range("E5").Value = 1.055713649400142
$ws.Range("F5").Value = 1.071610734183192
$ws.Range("I5").Value = 1.052817533617572
$ws.Range("J5").Value = 1.064274419007901
$ws.Range("K5").Value = 1.065155279148803
$ws.Range("L5").Value = 1.05809805050638
$ws.Range("M5").Value = 1.073957973561703
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06000608088409
$ws.Range("D6").Value = 1.062833311699996
$ws.Range("E6").Value = 1.055762892242888
$ws.Range("F6").Value = 1.071665557274681
$ws.Range("I6").Value = 1.05283741969458
$ws.Range("J6").Value = 1.064315410358249
$ws.Range("K6").Value = 1.065192458845241
$ws.Range("L6").Value = 1.058138650732802
$ws.Range("M6").Value = 1.074004283835238
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059614138195502
$ws.Range("D7").Value = 1.062518685529498
$ws.Range("E7").Value = 1.055424263760935
$ws.Range("F7").Value = 1.071288575799492
$ws.Range("I7").Value = 1.052700551533395
$ws.Range("J7").Value = 1.064033462872215
$ws.Range("K7").Value = 1.064936714512017
$ws.Range("L7").Value = 1.057859401942217
$ws.Range("M7").Value = 1.073685791139897
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.0579748696789
$ws.Range("D8").Value = 1.061202598985294
$ws.Range("E8").Value = 1.05400820783536
$ws.Range("F8").Value = 1.069712653702173
$ws.Range("I8").Value = 1.052125276030923
$ws.Range("J8").Value = 1.06285288352938
$ws.Range("K8").Value = 1.0638654997242
$ws.Range("L8").Value = 1.056690333990443
$ws.Range("M8").Value = 1.072353173805439
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.055082276726819
$ws.Range("D9").Value = 1.058879617230986
$ws.Range("E9").Value = 1.051510397225224
$ws.Range("F9").Value = 1.066934792788141
$ws.Range("I9").Value = 1.051099483141501
$ws.Range("J9").Value = 1.060764580479795
$ws.Range("K9").Value = 1.061969321715878
$ws.Range("L9").Value = 1.054623194613806
$ws.Range("M9").Value = 1.069999653091789
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053151028688239
$ws.Range("D10").Value = 1.057328267207135
$ws.Range("E10").Value = 1.049843376418913
$ws.Range("F10").Value = 1.065082169017708
$ws.Range("I10").Value = 1.050407436042756
$ws.Range("J10").Value = 1.059366911386689
$ws.Range("K10").Value = 1.060699354029712
$ws.Range("L10").Value = 1.053240232595083
$ws.Range("M10").Value = 1.068426979996791
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052314062191992
$ws.Range("D11").Value = 1.056655854459434
$ws.Range("E11").Value = 1.04912108336041
$ws.Range("F11").Value = 1.064279764697687
$ws.Range("I11").Value = 1.050105816672859
$ws.Range("J11").Value = 1.058760383944182
$ws.Range("K11").Value = 1.060148035291061
$ws.Range("L11").Value = 1.052640216845454
$ws.Range("M11").Value = 1.067745102018123
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.052003063514989
$ws.Range("D12").Value = 1.056405988105176
$ws.Range("E12").Value = 1.048852719939243
$ws.Range("F12").Value = 1.063981682561399
$ws.Range("I12").Value = 1.049993486310333
$ws.Range("J12").Value = 1.058534890994112
$ws.Range("K12").Value = 1.059943036539346
$ws.Range("L12").Value = 1.052417164311866
$ws.Range("M12").Value = 1.067491684669192
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052069778911977
$ws.Range("D13").Value = 1.056459589948283
$ws.Range("E13").Value = 1.048910288080421
$ws.Range("F13").Value = 1.064045623745112
$ws.Range("I13").Value = 1.050017594946793
$ws.Range("J13").Value = 1.058583269158117
$ws.Range("K13").Value = 1.059987019201783
$ws.Range("L13").Value = 1.052465018014842
$ws.Range("M13").Value = 1.067546049809799
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.052288357249788
$ws.Range("D14").Value = 1.056635202527371
$ws.Range("E14").Value = 1.04909890183007
$ws.Range("F14").Value = 1.064255125832042
$ws.Range("I14").Value = 1.050096537445498
$ws.Range("J14").Value = 1.058741748753042
$ws.Range("K14").Value = 1.060131094415708
$ws.Range("L14").Value = 1.052621782936775
$ws.Range("M14").Value = 1.067724157279264
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052423015628177
$ws.Range("D15").Value = 1.056743389615791
$ws.Range("E15").Value = 1.049215103430049
$ws.Range("F15").Value = 1.064384202475688
$ws.Range("I15").Value = 1.050145137329851
$ws.Range("J15").Value = 1.058839366474839
$ws.Range("K15").Value = 1.060219835438748
$ws.Range("L15").Value = 1.052718347069312
$ws.Range("M15").Value = 1.067833876880447
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.053206559869737
$ws.Range("D16").Value = 1.057372878785883
$ws.Range("E16").Value = 1.049891302707823
$ws.Range("F16").Value = 1.065135417396768
$ws.Range("I16").Value = 1.050427412157218
$ws.Range("J16").Value = 1.059407136458081
$ws.Range("K16").Value = 1.060735913274969
$ws.Range("L16").Value = 1.053280028540301
$ws.Range("M16").Value = 1.068472214833849
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.053697860426766
$ws.Range("D17").Value = 1.05776756049454
$ws.Range("E17").Value = 1.050315339393956
$ws.Range("F17").Value = 1.065606577757963
$ws.Range("I17").Value = 1.050603950481245
$ws.Range("J17").Value = 1.059762926599858
$ws.Range("K17").Value = 1.06105925524127
$ws.Range("L17").Value = 1.053632038014998
$ws.Range("M17").Value = 1.06887238477121
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053984358166449
$ws.Range("D18").Value = 1.057997707344724
$ws.Range("E18").Value = 1.05056262838806
$ws.Range("F18").Value = 1.065881377919235
$ws.Range("I18").Value = 1.050706733427319
$ws.Range("J18").Value = 1.059970325035991
$ws.Range("K18").Value = 1.06124771880737
$ws.Range("L18").Value = 1.053837245082009
$ws.Range("M18").Value = 1.069105710428633
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.05408203479444
$ws.Range("D19").Value = 1.058076170623507
$ws.Range("E19").Value = 1.05064694007434
$ws.Range("F19").Value = 1.065975074405156
$ws.Range("I19").Value = 1.050741747758022
$ws.Range("J19").Value = 1.060041020909584
$ws.Range("K19").Value = 1.061311956946805
$ws.Range("L19").Value = 1.053907196097113
$ws.Range("M19").Value = 1.069185253810992
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.053645155747605
$ws.Range("D20").Value = 1.057725221549868
$ws.Range("E20").Value = 1.050269848883253
$ws.Range("F20").Value = 1.065556028767725
$ws.Range("I20").Value = 1.05058502913547
$ws.Range("J20").Value = 1.059724766934607
$ws.Range("K20").Value = 1.061024577824758
$ws.Range("L20").Value = 1.05359428254804
$ws.Range("M20").Value = 1.068829459309019
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.05222399449766
$ws.Range("D21").Value = 1.05658349184026
$ws.Range("E21").Value = 1.049043361776632
$ws.Range("F21").Value = 1.064193433632469
$ws.Range("I21").Value = 1.050073299014656
$ws.Range("J21").Value = 1.058695086046714
$ws.Range("K21").Value = 1.060088673790931
$ws.Range("L21").Value = 1.052575624559241
$ws.Range("M21").Value = 1.067671712882853
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051329804021973
$ws.Range("D22").Value = 1.055865048080422
$ws.Range("E22").Value = 1.048271805191484
$ws.Range("F22").Value = 1.0633365212464
$ws.Range("I22").Value = 1.049749843679928
$ws.Range("J22").Value = 1.058046517006922
$ws.Range("K22").Value = 1.059498992220822
$ws.Range("L22").Value = 1.051934111576795
$ws.Range("M22").Value = 1.066942995051001
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.051803894054778
$ws.Range("D23").Value = 1.056245965578085
$ws.Range("E23").Value = 1.048680862054109
$ws.Range("F23").Value = 1.063790806010913
$ws.Range("I23").Value = 1.049921476014495
$ws.Range("J23").Value = 1.058390447160204
$ws.Range("K23").Value = 1.059811711911694
$ws.Range("L23").Value = 1.052274289237208
$ws.Range("M23").Value = 1.067329378439953
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053668970938334
$ws.Range("D24").Value = 1.057744352895255
$ws.Range("E24").Value = 1.050290404224228
$ws.Range("F24").Value = 1.065578869741457
$ws.Range("I24").Value = 1.050593579460739
$ws.Range("J24").Value = 1.05974201004005
$ws.Range("K24").Value = 1.06104024747785
$ws.Range("L24").Value = 1.053611342970953
$ws.Range("M24").Value = 1.068848855745557
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055830571778574
$ws.Range("D25").Value = 1.059480631872241
$ws.Range("E25").Value = 1.052156452988993
$ws.Range("F25").Value = 1.067653053767672
$ws.Range("I25").Value = 1.051366113430412
$ws.Range("J25").Value = 1.061305412576682
$ws.Range("K25").Value = 1.06246055325161
$ws.Range("L25").Value = 1.055158450506854
$ws.Range("M25").Value = 1.070608731483715
